# DOMA-4452: add "Meter place" column (S) for meter import example.
#
# The workbook has a single header row (row 1) describing import columns
# A..R, and 10 data rows (2..11) with sample meter data. This adds a new
# column S, "Meter place", with sample values "Kitchen" / "Bathroom"
# alternating down the data rows - matching the existing R ("Control
# readings date") column's header/data styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column R's formatting (fill/border/font per row) into column S so the
# new column matches the look of the rest of the imported table (header
# style in row 1, data style in rows 2-11), then overwrite the values.
$ws.Range("R1:R11").Copy() | Out-Null
$ws.Range("S1:S11").PasteSpecial(-4122, 0, $false, $false) | Out-Null

# New header
$ws.Range("S1").Value2 = "Meter place"

# Sample "meter place" values for the 10 sample data rows, alternating
# Kitchen / Bathroom.
$meterPlaces = @("Kitchen", "Bathroom", "Kitchen", "Bathroom", "Kitchen", "Bathroom", "Kitchen", "Bathroom", "Kitchen", "Bathroom")
for ($i = 0; $i -lt $meterPlaces.Length; $i++) {
    $row = 2 + $i
    $ws.Range("S$row").Value2 = $meterPlaces[$i]
}

# Match column S's width to column R's (new column should look the same
# width as the rest of the wide "date/place" columns).
$ws.Range("S1").ColumnWidth = $ws.Range("R1").ColumnWidth
